$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.700.17'
$ws.Range("E2").Value = '  +2.16%  '
$ws.Range("D3").Value = '1.806.68'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '''314.59'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").Value = '''0.9998'
$ws.Range("D7").Value = '''0.5402'
$ws.Range("E7").Value = '  -2.33%  '
$ws.Range("D8").Value = '''0.3794'
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("D9").Value = '''0.07544'
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").Value = '''42.70'
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("D11").Value = '''1.121'
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("E13").Value = '  -0.91%  '
$ws.Range("D14").Value = '''6.189'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").Value = '''7.393'
$ws.Range("E15").Value = '  +3.20%  '
$ws.Range("D16").Value = '1.796.06'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '''90.70'
$ws.Range("E17").Value = '  -1.45%  '
$ws.Range("E18").Value = '  -1.59%  '
$ws.Range("D19").Value = '''0.06447'
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").Value = '''0.9996'
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").Value = '''5.932'
$ws.Range("E22").Value = '  -0.68%  '
$ws.Range("D23").Value = '28.682.46'
$ws.Range("E23").Value = '  +2.00%  '
$ws.Range("D24").Value = '''11.22'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = '''2.104'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '''160.97'
$ws.Range("E26").Value = '  +3.08%  '
$ws.Range("D27").Value = '''20.54'
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("D28").Value = '''2.378'
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").Value = '2.007.55'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E31").Value = '  -3.39%  '
$ws.Range("D32").Value = '''0.1044'
$ws.Range("E32").Value = '  +0.39%  '
$ws.Range("D33").Value = '''5.690'
$ws.Range("E33").Value = '  -1.04%  '
$ws.Range("D34").Value = '''3.690'
$ws.Range("E34").Value = '  +2.58%  '
$ws.Range("D35").Value = '''0.2271'
$ws.Range("E35").Value = '  +6.61%  '
$ws.Range("D36").Value = '''0.06501'
$ws.Range("E36").Value = '  +7.59%  '
$ws.Range("D37").Value = '''8.984'
$ws.Range("E37").Value = '  +3.89%  '
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("D39").Value = '''5.062'
$ws.Range("E39").Value = '  +0.63%  '
$ws.Range("D40").Value = '''11.33'
$ws.Range("E40").Value = '  -2.03%  '
$ws.Range("D41").Value = '''0.6267'
$ws.Range("E41").Value = '  -0.72%  '
$ws.Range("D42").Value = '''1.205'
$ws.Range("E42").Value = '  +4.58%  '
$ws.Range("D43").Value = '''0.9992'
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("D44").Value = '''1.394'
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("D45").Value = '''13.35'
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("D46").Value = '''0.5894'
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("D47").Value = '''3.677'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("D48").Value = '''126.23'
$ws.Range("E48").Value = '  +3.40%  '
$ws.Range("D49").Value = '''1.965'
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("D50").Value = '''1.157'
$ws.Range("E50").Value = '  +1.71%  '
$ws.Range("E51").Value = '  +1.54%  '
